# Applies the "average/worst of ratio" summary block that was added to the
# freelancer k=0.2 distorted-greedy results sheet:
#   - J12           = AVERAGE(J2:J11)                         (bold)
#   - A14/B14       = "Average of SW(S*)/SW(OPT)" / AVERAGE(N2:N11)
#   - A15/B15       = "Average of SC(S*)/SC(OPT)" / AVERAGE(Z2:Z11)
#   - A16/B16       = "Worst of SW(S*)/SW(OPT)"   / MIN(N2:N11)
#   - A17/B17       = "Worst of SC(S*)/SC(OPT)"   / MAX(Z2:Z11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: first summary label/value pair - this one defines the style --
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$b14 = $ws.Range("B14")
$b14.Formula = "=AVERAGE(N2:N11)"
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108
$ws.Rows.Item(14).RowHeight = 15.6

# --- Row 12: overall average of the k fraction column (J) -----------------
$avgK = $ws.Range("J12")
$avgK.Formula = "=AVERAGE(J2:J11)"
$avgK.Font.Bold = $true

# --- Rows 15-17: remaining labels/values, reusing row 14's formatting -----
$b14.Copy() | Out-Null

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$b15 = $ws.Range("B15")
$b15.Formula = "=AVERAGE(Z2:Z11)"
$b15.PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(15).RowHeight = 15.6

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$b16 = $ws.Range("B16")
$b16.Formula = "=MIN(N2:N11)"
$b16.PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(16).RowHeight = 15.6

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$b17 = $ws.Range("B17")
$b17.Formula = "=MAX(Z2:Z11)"
$b17.PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(17).RowHeight = 15.6

$excel.CutCopyMode = $false

# --- Page setup used for the resave (portrait, paper size 9 = A4) ---------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection left on the newly added summary block -----------------------
$ws.Range("A14:B17").Select() | Out-Null
